$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

$ws.Range("B11").Value = 0.1178764178612974

$ws.Range("B12").Value = 0.3408936812763721
$ws.Range("C12").Value = "{'codebleu': 0.3408936812763721, 'ngram_match_score': 0.11671135607664002, 'weighted_ngram_match_score': 0.15620402836950756, 'syntax_match_score': 0.554945054945055, 'dataflow_match_score': 0.5357142857142857}"

$ws.Range("B13").Value = 0.8509027962175195
